$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "159617931477214"
$ws.Range("D2").Value = "762027117282800"
$ws.Range("C3").Value = "1761231377490564"
$ws.Range("D3").Value = "616734738469807"
$ws.Range("C4").Value = "264846833886893"
$ws.Range("D4").Value = "1544909859170915"
$ws.Range("C5").Value = "1549979828593724"
$ws.Range("D5").Value = "332252896899497"
$ws.Range("C6").Value = "1455330608067225"
$ws.Range("D6").Value = "1856620574550792"
$ws.Range("C7").Value = "316603388502163"
$ws.Range("D7").Value = "1665045370384962"
$ws.Range("C8").Value = "784012071764293"
$ws.Range("D8").Value = "159617931477214"
$ws.Range("C9").Value = "1669243136660143"
$ws.Range("D9").Value = "2116731221905896"
$ws.Range("C10").Value = "254776108322246"
$ws.Range("D10").Value = "784012071764293"
$ws.Range("C11").Value = "1798895427089724"
$ws.Range("D11").Value = "206240536583250"
$ws.Range("C12").Value = "774067269282955"
$ws.Range("D12").Value = "1384506521580724"
$ws.Range("C13").Value = "362228173958804"
$ws.Range("D13").Value = "494023894278277"
$ws.Range("C14").Value = "107407153282844"
$ws.Range("D14").Value = "737452329661021"
$ws.Range("C15").Value = "585539334989100"
$ws.Range("D15").Value = "1720869341538398"
$ws.Range("C16").Value = "1208088295924736"
$ws.Range("D16").Value = "322008967984765"
$ws.Range("C17").Value = "391037157727116"
$ws.Range("D17").Value = "277314789020591"
$ws.Range("C18").Value = "179112262432141"
$ws.Range("D18").Value = "1798895427089724"
$ws.Range("C19").Value = "703894706306308"
$ws.Range("D19").Value = "1208088295924736"
$ws.Range("C20").Value = "492401487499641"
$ws.Range("D20").Value = "316603388502163"
$ws.Range("C21").Value = "1544909859170915"
$ws.Range("D21").Value = "249978345123454"
$ws.Range("C22").Value = "625260624533249"
$ws.Range("D22").Value = "436038656468443"
$ws.Range("C23").Value = "1384506521580724"
$ws.Range("D23").Value = "173960829412142"
$ws.Range("C24").Value = "285153488580341"
$ws.Range("D24").Value = "1362279090456478"
$ws.Range("C25").Value = "459470627407885"
$ws.Range("D25").Value = "492401487499641"
$ws.Range("C26").Value = "415793858487646"
$ws.Range("D26").Value = "264846833886893"
$ws.Range("C27").Value = "761161927309184"
$ws.Range("D27").Value = "246933925685884"
$ws.Range("C28").Value = "711816275694043"
$ws.Range("D28").Value = "275712962530384"
$ws.Range("C29").Value = "145971692474202"
$ws.Range("D29").Value = "179112262432141"
$ws.Range("C30").Value = "307963623358864"
$ws.Range("D30").Value = "264292347041953"
$ws.Range("C31").Value = "1792598261018671"
$ws.Range("D31").Value = "781074362021860"
$ws.Range("C32").Value = "319718041464688"
$ws.Range("D32").Value = "319718041464688"
$ws.Range("C33").Value = "1741096309512211"
$ws.Range("D33").Value = "108757516583030"
$ws.Range("C34").Value = "876689279025944"
$ws.Range("D34").Value = "1911119432502301"
$ws.Range("C35").Value = "451122651649814"
$ws.Range("D35").Value = "1471037103142820"
$ws.Range("C36").Value = "498723016852762"
$ws.Range("D36").Value = "498723016852762"
$ws.Range("C37").Value = "322008967984765"
$ws.Range("D37").Value = "451122651649814"
$ws.Range("C38").Value = "141059952712333"
$ws.Range("D38").Value = "1455330608067225"
$ws.Range("C39").Value = "423158777729623"
$ws.Range("D39").Value = "459470627407885"
$ws.Range("C40").Value = "603879913082562"
$ws.Range("D40").Value = "171076470160303"
$ws.Range("C41").Value = "386412421458390"
$ws.Range("D41").Value = "154362801673336"
$ws.Range("C42").Value = "1471037103142820"
$ws.Range("D42").Value = "703894706306308"
$ws.Range("C43").Value = "263245977200582"
$ws.Range("D43").Value = "1826787617399125"
$ws.Range("C44").Value = "1665045370384962"
$ws.Range("D44").Value = "1209295105801230"
$ws.Range("C45").Value = "781074362021860"
$ws.Range("D45").Value = "126348794234977"
$ws.Range("C46").Value = "108757516583030"
$ws.Range("D46").Value = "1761231377490564"
$ws.Range("C47").Value = "154362801673336"
$ws.Range("D47").Value = "145971692474202"
$ws.Range("C48").Value = "246933925685884"
$ws.Range("D48").Value = "205882209595847"
$ws.Range("C49").Value = "206240536583250"
$ws.Range("D49").Value = "761161927309184"
$ws.Range("C50").Value = "264292347041953"
$ws.Range("D50").Value = "423158777729623"
$ws.Range("C51").Value = "737452329661021"
$ws.Range("D51").Value = "774067269282955"
$ws.Range("C52").Value = "383453621698654"
$ws.Range("D52").Value = "1549979828593724"
$ws.Range("C53").Value = "174548946759532"
$ws.Range("D53").Value = "415793858487646"
$ws.Range("C54").Value = "1720869341538398"
$ws.Range("D54").Value = "386412421458390"
$ws.Range("C55").Value = "173960829412142"
$ws.Range("C56").Value = "922749794537289"
$ws.Range("C57").Value = "184085025268369"
$ws.Range("C58").Value = "332252896899497"
$ws.Range("C59").Value = "383831868477620"
$ws.Range("C60").Value = "2116731221905896"
$ws.Range("C61").Value = "2153529108017980"
$ws.Range("C62").Value = "1826787617399125"
$ws.Range("C63").Value = "218683288246409"
$ws.Range("C64").Value = "517772401618874"

[void]$ws.Range("G6").Select()
